# Generate Report for Handback
# Updates the handback-status report with refreshed generation timestamps
# and a corrected translation-type flag ("ht" -> "mt") for the
# 69d53e7a-aebb-4003-9e0c-7d3fea20b1c4 row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" column (G) for the
# 69d53e7a... row (row 3) and the 7bdf088b... row (row 4), which shared
# the same timestamp text.
$wsOverview.Range("G3").Value = "2016-08-22 22:16:03"
$wsOverview.Range("G4").Value = "2016-08-22 22:16:03"

# zh-cn sheet, row 3 & 4 (69d53e7a... / 7bdf088b...)
# Priority column (E): "ht" -> "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H3").Value = "2016-08-22 22:15:56"
$wsZhCn.Range("H4").Value = "2016-08-22 22:15:56"
# Correspond Handback DateTime column (K)
$wsZhCn.Range("K3").Value = "2016-08-22 22:16:28"
$wsZhCn.Range("K4").Value = "2016-08-22 22:16:28"

# de-de sheet, row 3 & 4 (69d53e7a... / 7bdf088b...)
# Priority column (E): "ht" -> "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsDeDe.Range("H3").Value = "2016-08-22 22:16:03"
$wsDeDe.Range("H4").Value = "2016-08-22 22:16:03"
# Correspond Handback DateTime column (K)
$wsDeDe.Range("K3").Value = "2016-08-22 22:16:35"
$wsDeDe.Range("K4").Value = "2016-08-22 22:16:35"
